$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.205438613891602
$ws.Range("B1").Value = 2.617996692657471
$ws.Range("D1").Value = 2.171185970306396
$ws.Range("E1").Value = 1.167295694351196
